# Ponoko case price increased ($48 -> $42.16 net of the 38% Ponoko discount on
# a new $68 list price); reflect it in both the summary sheet and the detailed
# cost breakdown sheet, and leave parts_1 as the active/selected sheet.
$wb = $excel.ActiveWorkbook

$wsParts1  = $wb.Worksheets.Item("parts_1")
$wsParts10 = $wb.Worksheets.Item("parts_10")

# --- parts_10 ("Wood case" cost block) ---
$wsParts10.Activate()

# New price breakdown rows: raw Ponoko price, then price after their 38% off.
$wsParts10.Range("B24").Value = 68
$wsParts10.Range("B25").Formula = "=B24-(B24*38%)"

# Update the "Notes" cell that documents the Ponoko bulk price.
$wsParts10.Range("G4").Value = "Or order 10+ from Ponoko for @`$42.16"

$wsParts10.Range("G5").Select()

# --- parts_1 ("Case" row): bump the case cost to match the new Ponoko price ---
$wsParts1.Range("C2").Value = 68

$wsParts1.Activate()
$wsParts1.Range("C3").Select()
